$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 199.33333
$ws.Range("I4").Value = 199.33333
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 199.33333
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -85.33332999999999
$ws.Range("N4").ClearContents()

$ws.Range("H33").Value = 135.88889
$ws.Range("I33").Value = 135.88889
$ws.Range("K33").Value = 135.88889
$ws.Range("M33").Value = 93.11111

$ws.Range("H53").Value = 355.4
$ws.Range("I53").Value = 332.76923
$ws.Range("J53").Value = 502.5
$ws.Range("K53").Value = 332.76923
$ws.Range("L53").Value = 502.5
$ws.Range("M53").Value = 304.23077
$ws.Range("N53").Value = -1776.5

$ws.Range("H70").Value = 998
$ws.Range("I70").Value = 998
$ws.Range("K70").Value = 2994
$ws.Range("M70").Value = -2724

$ws.Range("H73").Value = 998
$ws.Range("I73").Value = 998
$ws.Range("K73").Value = 2994
$ws.Range("M73").Value = -2058

$ws.Range("H100").Value = 1846.4
$ws.Range("J100").Value = 1497.7142
$ws.Range("L100").Value = 1497.7142
$ws.Range("N100").Value = -2579.7142

$ws.Range("H107").Value = 290.42856
$ws.Range("I107").Value = 305.66666
$ws.Range("K107").Value = 305.66666
$ws.Range("M107").Value = 1614.33334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 739.9375
$ws.Range("I2").Value = 804.6667
$ws.Range("J2").Value = 545.75
$ws.Range("K2").Value = 804.6667
$ws.Range("L2").Value = 545.75
$ws.Range("M2").Value = -691.6667
$ws.Range("N2").Value = -771.75

$ws.Range("H4").Value = 915
$ws.Range("I4").Value = 1092
$ws.Range("J4").Value = 797
$ws.Range("K4").Value = 1092
$ws.Range("L4").Value = 797
$ws.Range("M4").Value = -976
$ws.Range("N4").Value = -1029

$ws.Range("H32").Value = 7045.8
$ws.Range("I32").Value = 5374.3228
$ws.Range("K32").Value = 5374.3228
$ws.Range("M32").Value = -5087.3228

$ws.Range("H41").Value = 800
$ws.Range("I41").Value = 800
$ws.Range("K41").Value = 800
$ws.Range("M41").Value = -386

$ws.Range("H61").Value = 2996.8
$ws.Range("I61").Value = 2328.8333
$ws.Range("K61").Value = 2328.8333
$ws.Range("M61").Value = -2116.8333

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H74").Value = 39987940
$ws.Range("I74").Value = 66642230
$ws.Range("K74").Value = 66642230
$ws.Range("M74").Value = -66641356

$ws.Range("H77").Value = 39987940
$ws.Range("I77").Value = 66642230
$ws.Range("K77").Value = 333211150
$ws.Range("M77").Value = -333206782

$ws.Range("H110").Value = 1084.2
$ws.Range("I110").Value = 1084.2
$ws.Range("K110").Value = 1084.2
$ws.Range("M110").Value = 960.8

$ws.Range("H116").Value = 739.9375
$ws.Range("I116").Value = 804.6667
$ws.Range("J116").Value = 545.75
$ws.Range("K116").Value = 804.6667
$ws.Range("L116").Value = 545.75
$ws.Range("M116").Value = 1489.3333
$ws.Range("N116").Value = -5133.75

$ws.Range("H132").Value = 2497.5557
$ws.Range("I132").Value = 1818.2778
$ws.Range("J132").Value = 3856.111
$ws.Range("K132").Value = 5454.8334
$ws.Range("L132").Value = 11568.333
$ws.Range("M132").Value = -2924.8334
$ws.Range("N132").Value = -16628.333

$ws.Range("H136").Value = 2996.8
$ws.Range("I136").Value = 2328.8333
$ws.Range("K136").Value = 6986.499899999999
$ws.Range("M136").Value = -4436.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 739.9375
$ws.Range("I3").Value = 804.6667
$ws.Range("J3").Value = 545.75
$ws.Range("K3").Value = 804.6667
$ws.Range("L3").Value = 545.75
$ws.Range("M3").Value = -690.6667
$ws.Range("N3").Value = -773.75

$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 100
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 100
$ws.Range("L8").Value = 100
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = -380

$ws.Range("H11").Value = 164
$ws.Range("I11").Value = 164
$ws.Range("K11").Value = 164
$ws.Range("M11").Value = -24

$ws.Range("H86").Value = 3924.7
$ws.Range("I86").Value = 2405.875
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 2405.875
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -1282.875
$ws.Range("N86").Value = -12246

$ws.Range("H89").Value = 3924.7
$ws.Range("I89").Value = 2405.875
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 12029.375
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -6413.375
$ws.Range("N89").Value = -61232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 468.66666
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1400
$ws.Range("M22").Value = 347
$ws.Range("N22").Value = -2100

$ws.Range("H31").Value = 1111
$ws.Range("I31").Value = 1111
$ws.Range("K31").Value = 1111
$ws.Range("M31").Value = -816

$ws.Range("H34").Value = 1111
$ws.Range("I34").Value = 1111
$ws.Range("K34").Value = 1111
$ws.Range("M34").Value = -909

$ws.Range("H107").Value = 1792.1818
$ws.Range("I107").Value = 980.7143
$ws.Range("K107").Value = 980.7143
$ws.Range("M107").Value = 939.2857

$ws.Range("H109").Value = 61283.5
$ws.Range("J109").Value = 61283.5
$ws.Range("L109").Value = 61283.5
$ws.Range("N109").Value = -63363.5

$ws.Range("H134").Value = 3096.6667
$ws.Range("I134").Value = 2981.3845
$ws.Range("J134").Value = 3396.4
$ws.Range("K134").Value = 8944.1535
$ws.Range("L134").Value = 10189.2
$ws.Range("M134").Value = -6409.1535
$ws.Range("N134").Value = -15259.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 163
$ws.Range("I2").Value = 194.5
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 1167
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -1054
$ws.Range("N2").Value = -826

$ws.Range("H34").Value = 3199.8
$ws.Range("I34").Value = 1566.3334
$ws.Range("K34").Value = 4699.0002
$ws.Range("M34").Value = -4615.0002

$ws.Range("H51").Value = 1499.5
$ws.Range("I51").Value = 1499.5
$ws.Range("K51").Value = 4498.5
$ws.Range("M51").Value = -4038.5

$ws.Range("H140").Value = 1060.5555
$ws.Range("I140").Value = 1060.5555
$ws.Range("K140").Value = 3181.6665
$ws.Range("M140").Value = 1998.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2307.5217
$ws.Range("I132").Value = 1475.1177
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 4425.3531
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -1895.3531
$ws.Range("N132").Value = -19058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 99000
$ws.Range("J43").Value = 99000
$ws.Range("L43").Value = 99000
$ws.Range("N43").Value = -99386

$ws.Range("H132").Value = 3630.6428
$ws.Range("I132").Value = 3003.6667
$ws.Range("K132").Value = 9011.000100000001
$ws.Range("M132").Value = -6481.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 60299.5
$ws.Range("J27").Value = 60299.5
$ws.Range("L27").Value = 60299.5
$ws.Range("N27").Value = -60437.5

$ws.Range("H107").Value = 499.5
$ws.Range("I107").Value = 499.5
$ws.Range("K107").Value = 1498.5
$ws.Range("M107").Value = 421.5
